$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, pushing existing rows 16..108 down to 17..109.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new record's data.
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44881
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 100112021
$ws.Range("G16").Value = "Ají"
$ws.Range("H16").Value = "Inferno"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("N16").Value = "$/caja 10 kilos"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 1800
$ws.Range("Q16").Value = 10
$ws.Range("R16").Value = "Hortaliza"
